{"js": "// Apply 4 targeted text replacements inside DOCPROPERTY field-result runs.\n// Each search string is unique in the document, so a plain body.search()\n// with exact matching is sufficient and safe.\n\nconst replacements = [\n  { find: \"01 de noviembre de 2017.\", replace: \"12 de noviembre de 2017.\" },\n  {\n    find: \"Fiscalia Nacional en lo Criminal y Correccional Nro 1\",\n    replace: \"Fiscalia Nacional en lo Criminal y Correccional Nro 3\",\n  },\n  { find: \"$ 1.500,00\", replace: \"$ 43.000,00\" },\n  {\n    find:\n      \"Finalmente, la presente erogaci\u00f3n de fondos es solicitada por este curso debido a que rrr\",\n    replace: \" \",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply 4 targeted text replacements inside DOCPROPERTY field-result runs.\n# Each Find pattern is unique in the document, so plain Find/Replace\n# (MatchCase, no wildcards) is sufficient and safe.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0          # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        [ref]$findText,        # FindText\n        [ref]$true,            # MatchCase\n        [ref]$false,           # MatchWholeWord\n        [ref]$false,           # MatchWildcards\n        [ref]$false,           # MatchSoundsLike\n        [ref]$false,           # MatchAllWordForms\n        [ref]$true,            # Forward\n        [ref]0,                # Wrap = wdFindStop\n        [ref]$false,           # Format\n        [ref]$replaceText,     # ReplaceWith\n        [ref]2                 # Replace = wdReplaceAll\n    ) | Out-Null\n\n    if (-not $find.Found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-ExactText '01 de noviembre de 2017.' '12 de noviembre de 2017.'\nReplace-ExactText 'Fiscalia Nacional en lo Criminal y Correccional Nro 1' 'Fiscalia Nacional en lo Criminal y Correccional Nro 3'\nReplace-ExactText '$ 1.500,00' '$ 43.000,00'\nReplace-ExactText 'Finalmente, la presente erogaci\u00f3n de fondos es solicitada por este curso debido a que rrr' ' '\n"}
